$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates: Volume/Number and report week date range ---
$ws.Range("A8").Value = "Volume 32   Number  9"
$ws.Range("C9").Value = "Report Covering the Week  2/24/2025  Through  3/2/2025"

# --- Crime Complaints data table (rows 14-31) ---
# Stable style donor cells (never themselves edited by this script):
#   C15 = style 13 (text / "N/A" placeholder style)
#   F15 = style 14 (integer count style)
#   L15 = style 15 (percent-change style)

$ws.Range("D14").Value = "'0"
$ws.Range("C15").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null
$ws.Range("E14").Value = "***.*"
$ws.Range("C15").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null
$ws.Range("D15").Value = 1
$ws.Range("F15").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("E15").Value = -100
$ws.Range("L15").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = 0
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = -7.142857142857
$ws.Range("I16").Value = 26
$ws.Range("J16").Value = 34
$ws.Range("K16").Value = -23.529411764705
$ws.Range("L16").Value = -16.129032258064
$ws.Range("M16").Value = 13.043478260869
$ws.Range("N16").Value = -87.735849056603
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 116.666666666667
$ws.Range("I17").Value = 25
$ws.Range("J17").Value = 12
$ws.Range("K17").Value = 108.333333333333
$ws.Range("L17").Value = -21.875
$ws.Range("M17").Value = 127.272727272727
$ws.Range("N17").Value = -30.555555555555
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -71.428571428571
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 28
$ws.Range("H18").Value = -46.428571428571
$ws.Range("I18").Value = 43
$ws.Range("J18").Value = 47
$ws.Range("K18").Value = -8.510638297872
$ws.Range("L18").Value = -6.521739130434
$ws.Range("M18").Value = -31.746031746031
$ws.Range("N18").Value = -90.611353711790
$ws.Range("C19").Value = 33
$ws.Range("D19").Value = 32
$ws.Range("E19").Value = 3.125
$ws.Range("F19").Value = 96
$ws.Range("G19").Value = 110
$ws.Range("H19").Value = -12.727272727272
$ws.Range("I19").Value = 225
$ws.Range("J19").Value = 265
$ws.Range("K19").Value = -15.094339622641
$ws.Range("L19").Value = -6.639004149377
$ws.Range("M19").Value = 15.384615384615
$ws.Range("N19").Value = -64.171974522293
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = -33.333333333333
$ws.Range("I20").Value = 9
$ws.Range("J20").Value = 9
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = -55
$ws.Range("M20").Value = 50
$ws.Range("N20").Value = -98.553054662379
$ws.Range("C21").Value = 44
$ws.Range("D21").Value = 47
$ws.Range("E21").Value = -6.382978723404
$ws.Range("F21").Value = 142
$ws.Range("G21").Value = 167
$ws.Range("H21").Value = -14.970059880239
$ws.Range("I21").Value = 331
$ws.Range("J21").Value = 371
$ws.Range("K21").Value = -10.781671159029
$ws.Range("L21").Value = -11.021505376344
$ws.Range("M21").Value = 9.966777408637
$ws.Range("N21").Value = -83.120856705762
$ws.Range("C22").Value = "'0"
$ws.Range("C15").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("D22").Value = "'0"
$ws.Range("C15").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").Value = "***.*"
$ws.Range("C15").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 66.666666666666
$ws.Range("L22").Value = 14.285714285714
$ws.Range("D23").Value = 2
$ws.Range("F15").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4122) | Out-Null
$ws.Range("E23").Value = -50
$ws.Range("L15").Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4122) | Out-Null
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -50
$ws.Range("I23").Value = 4
$ws.Range("J23").Value = 6
$ws.Range("K23").Value = -33.333333333333
$ws.Range("L23").Value = 100
$ws.Range("M23").Value = -20
$ws.Range("C24").Value = 44
$ws.Range("D24").Value = 57
$ws.Range("E24").Value = -22.807017543859
$ws.Range("F24").Value = 223
$ws.Range("G24").Value = 185
$ws.Range("H24").Value = 20.540540540540
$ws.Range("I24").Value = 588
$ws.Range("J24").Value = 484
$ws.Range("K24").Value = 21.487603305785
$ws.Range("L24").Value = 4.626334519572
$ws.Range("M24").Value = 133.333333333333
$ws.Range("C25").Value = 46
$ws.Range("D25").Value = 55
$ws.Range("E25").Value = -16.363636363636
$ws.Range("F25").Value = 201
$ws.Range("G25").Value = 177
$ws.Range("H25").Value = 13.559322033898
$ws.Range("I25").Value = 496
$ws.Range("J25").Value = 449
$ws.Range("K25").Value = 10.467706013363
$ws.Range("L25").Value = 5.084745762711
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -12.5
$ws.Range("F26").Value = 23
$ws.Range("G26").Value = 27
$ws.Range("H26").Value = -14.814814814814
$ws.Range("I26").Value = 41
$ws.Range("J26").Value = 56
$ws.Range("K26").Value = -26.785714285714
$ws.Range("L26").Value = -31.666666666666
$ws.Range("M26").Value = -32.786885245901
$ws.Range("D27").Value = 1
$ws.Range("F15").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("E27").Value = -100
$ws.Range("L15").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 3
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = -40
$ws.Range("C28").Value = "'0"
$ws.Range("C15").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 12
$ws.Range("H28").Value = -41.666666666666
$ws.Range("J28").Value = 21
$ws.Range("K28").Value = -23.809523809523
$ws.Range("D29").Value = 1
$ws.Range("F15").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null
$ws.Range("E29").Value = -100
$ws.Range("L15").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4122) | Out-Null
$ws.Range("G29").Value = 1
$ws.Range("F15").Copy() | Out-Null
$ws.Range("G29").PasteSpecial(-4122) | Out-Null
$ws.Range("H29").Value = -100
$ws.Range("L15").Copy() | Out-Null
$ws.Range("H29").PasteSpecial(-4122) | Out-Null
$ws.Range("J29").Value = 1
$ws.Range("F15").Copy() | Out-Null
$ws.Range("J29").PasteSpecial(-4122) | Out-Null
$ws.Range("K29").Value = -100
$ws.Range("L15").Copy() | Out-Null
$ws.Range("K29").PasteSpecial(-4122) | Out-Null
$ws.Range("D30").Value = 1
$ws.Range("F15").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null
$ws.Range("E30").Value = -100
$ws.Range("L15").Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4122) | Out-Null
$ws.Range("G30").Value = 1
$ws.Range("F15").Copy() | Out-Null
$ws.Range("G30").PasteSpecial(-4122) | Out-Null
$ws.Range("H30").Value = -100
$ws.Range("L15").Copy() | Out-Null
$ws.Range("H30").PasteSpecial(-4122) | Out-Null
$ws.Range("J30").Value = 1
$ws.Range("F15").Copy() | Out-Null
$ws.Range("J30").PasteSpecial(-4122) | Out-Null
$ws.Range("K30").Value = -100
$ws.Range("L15").Copy() | Out-Null
$ws.Range("K30").PasteSpecial(-4122) | Out-Null
$ws.Range("D31").Value = 1
$ws.Range("F15").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4122) | Out-Null
$ws.Range("E31").Value = -100
$ws.Range("L15").Copy() | Out-Null
$ws.Range("E31").PasteSpecial(-4122) | Out-Null
$ws.Range("F31").Value = 1
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 4
$ws.Range("K31").Value = -50

Write-Output "Applied weekly CompStat update."
